# Agregada tabla independiente para referentes
#
# The original "Hoja1" sheet mixed two different entities in one table:
# the artwork ("obra") columns (A-E) and the newspaper-clipping
# ("referente") columns (F-I). This splits them into two proper sheets:
#   - "obras"      : Título, Fecha, Dimensiones, Técnica, Archivo
#   - "referentes" : Título, Fecha, Periódico, Archivo

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "obras"

# --- Clear out the old "referente" columns (F:I) from the obras sheet ---
# Rows 2-4 lose their F:I content entirely (cells disappear).
$ws1.Range("F2:I4").ClearContents()
# Row 1 (header) keeps an empty, but still-bold-styled, F1:I1 cell.
$ws1.Range("F1:I1").Font.Bold = $true
$ws1.Range("F1:I1").ClearContents()

# Bold the remaining obras header row A1:E1
$ws1.Range("A1:E1").Font.Bold = $true

# --- Add the new "referentes" sheet right after "obras" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "referentes"
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws2.Columns.Item(1).ColumnWidth = 54
$ws2.Columns.Item(2).ColumnWidth = 17.66666667
$ws2.Columns.Item(3).ColumnWidth = 16.83333333
$ws2.Columns.Item(4).ColumnWidth = 33.33333333

$ws2.Range("A1").Value = "Título"
$ws2.Range("B1").Value = "Fecha"
$ws2.Range("C1").Value = "Periódico"
$ws2.Range("D1").Value = "Archivo"
$ws2.Range("A1:D1").Font.Bold = $true

$ws2.Range("A2").Value = "Doble suicidio en ""El Sisga"""
$ws2.Range("B2").Value = "Junio 29 1965"
$ws2.Range("C2").Value = "El Tiempo"
$ws2.Range("D2").Value = "doble-suicidio-el-tiempo.jpg"

$ws2.Range("A3").Value = "Una indígena y su hijo murieron en persecución"
$ws2.Range("B3").Value = "Mayo 24 del 96"
$ws2.Range("C3").Value = "El Tiempo"
$ws2.Range("D3").Value = "indigena-hijo-el-tiempo.jpg"

$ws2.Range("A4").Value = " Láminas de paisajes latinoamericanos"
$ws2.Range("D4").Value = " laminas-paisajes.jpg"

$ws2.Range("A5").Value = " Exmilitar Mata a la Ezposa de su Amigo y se Suicida"
$ws2.Range("D5").Value = "exmilitar-mata-esposa.jpg"

# --- Selection / active-sheet state to match the authored workbook ---
$ws1.Range("F13").Select()
$ws2.Range("C12").Select()
